$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9894965887069702
$ws.Range("B1").Value = 1.521806001663208
$ws.Range("D1").Value = 1.754169821739197
$ws.Range("E1").Value = 1.049919247627258
